# Swap the B/D/E/F/G column values between paired rows in the stock report.
# Each pair of rows represents the same product (same description in column C)
# but the batch number (B), rate (D), MRP (E), quantity (F) and value (G)
# were recorded against the wrong row and need to be swapped back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(149, 150),
    @(183, 184),
    @(279, 280),
    @(316, 318),
    @(346, 347),
    @(350, 352),
    @(355, 356),
    @(372, 373),
    @(379, 380),
    @(382, 383),
    @(389, 390),
    @(419, 420),
    @(421, 422),
    @(581, 582),
    @(590, 591),
    @(601, 602),
    @(687, 688)
)

$cols = @("B", "D", "E", "F", "G")

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()

        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
